$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"1"
$ws.Range("G2").Value = [double]"0.2757802710344801"
$ws.Range("J2").Value = [double]"0.2476798937108692"
$ws.Range("S2").Value = [double]"-2.578523513553466"
$ws.Range("U2").Value = [double]"2.578523513553466"
$ws.Range("C3").Value = [double]"47.10738243223267"
$ws.Range("E3").Value = [double]"1"
$ws.Range("G3").Value = [double]"0.3431638519279241"
$ws.Range("J3").Value = [double]"0.2881100616251468"
$ws.Range("S3").Value = [double]"-2.720163665405287"
$ws.Range("U3").Value = [double]"2.720163665405287"
$ws.Range("C4").Value = [double]"33.50656410520624"
$ws.Range("E4").Value = [double]"1"
$ws.Range("G4").Value = [double]"0.3962168548245172"
$ws.Range("J4").Value = [double]"0.33018120641142"
$ws.Range("S4").Value = [double]"-2.701311175357899"
$ws.Range("U4").Value = [double]"2.701311175357899"
$ws.Range("C5").Value = [double]"20.00000822841674"
$ws.Range("D5").Value = [double]"-8.836078093497032E-08"
$ws.Range("E5").Value = [double]"1"
$ws.Range("F5").Value = [double]"4.418039046748516E-09"
$ws.Range("G5").Value = [double]"-0.291145193183945"
$ws.Range("H5").Value = [double]"7.9285868495038E-09"
$ws.Range("J5").Value = [double]"1.397380953274697E-07"
$ws.Range("R5").Value = [double]"8.836078093497032E-08"
$ws.Range("S5").Value = [double]"-1.557322567548503E-06"
$ws.Range("U5").Value = [double]"1.557322567548503E-06"
$ws.Range("E6").Value = [double]"1"
$ws.Range("E7").Value = [double]"1"
$ws.Range("G7").Value = [double]"0.4163950387989465"
$ws.Range("J7").Value = [double]"0.3469959791983928"
$ws.Range("S7").Value = [double]"-3.351497329532523"
$ws.Range("U7").Value = [double]"3.352117109752303"
$ws.Range("C8").Value = [double]"51.09695633863818"
$ws.Range("D8").Value = [double]"-2.371889295868397"
$ws.Range("E8").Value = [double]"1"
$ws.Range("F8").Value = [double]"0.1185944647934199"
$ws.Range("G8").Value = [double]"0.8100065000603751"
$ws.Range("H8").Value = [double]"0.2880797662138852"
$ws.Range("J8").Value = [double]"0.4910670858704626"
$ws.Range("R8").Value = [double]"2.371889295868397"
$ws.Range("S8").Value = [double]"-3.847501971859238"
$ws.Range("U8").Value = [double]"4.043174499331766"
$ws.Range("E9").Value = [double]"1"
$ws.Range("G9").Value = [double]"-1.12250361326851"
$ws.Range("J9").Value = [double]"0.1174132855466139"
$ws.Range("O9").Value = [double]"2.572754361107599E-06"
$ws.Range("P9").Value = [double]"9.570886024505796"
$ws.Range("T9").Value = [double]"-1.406064460212672"
$ws.Range("U9").Value = [double]"1.406064460212672"
$ws.Range("E10").Value = [double]"1"
$ws.Range("G10").Value = [double]"-0.2994219358709864"
$ws.Range("J10").Value = [double]"0.2440891495234913"
$ws.Range("O10").Value = [double]"7.093727401041861E-07"
$ws.Range("P10").Value = [double]"4.429110693367102"
$ws.Range("T10").Value = [double]"-2.921491598319568"
$ws.Range("U10").Value = [double]"2.921491598319568"
$ws.Range("E11").Value = [double]"1"
$ws.Range("G11").Value = [double]"0.5833178818634276"
$ws.Range("J11").Value = [double]"0.4860983375149515"
$ws.Range("S11").Value = [double]"-0.4053217716707529"
$ws.Range("U11").Value = [double]"4.518039354088335"
$ws.Range("C12").Value = [double]"87.97339114164625"
$ws.Range("E12").Value = [double]"1"
$ws.Range("G12").Value = [double]"0.7362097884893807"
$ws.Range("J12").Value = [double]"0.6135083106679333"
$ws.Range("S12").Value = [double]"-0.8290554145674749"
$ws.Range("U12").Value = [double]"5.503897172809233"
$ws.Range("C13").Value = [double]"83.82811406880887"
$ws.Range("D13").Value = [double]"-4.140548588477222"
$ws.Range("E13").Value = [double]"1"
$ws.Range("F13").Value = [double]"0.2070274294238611"
$ws.Range("G13").Value = [double]"1.600573204122853"
$ws.Range("H13").Value = [double]"0.5797785927913192"
$ws.Range("J13").Value = [double]"0.8506623249894878"
$ws.Range("R13").Value = [double]"4.140548588477222"
$ws.Range("S13").Value = [double]"-1.072204774671102"
$ws.Range("U13").Value = [double]"6.075092686759014"
$ws.Range("C14").Value = [double]"57.76434725306725"
$ws.Range("E14").Value = [double]"1"
$ws.Range("G14").Value = [double]"1.031830711928773"
$ws.Range("J14").Value = [double]"0.8598590174739729"
$ws.Range("S14").Value = [double]"-1.726084137496198"
$ws.Range("U14").Value = [double]"6.424864357276418"
$ws.Range("C15").Value = [double]"49.13392656558626"
$ws.Range("E15").Value = [double]"1"
$ws.Range("G15").Value = [double]"0.7306121673939387"
$ws.Range("J15").Value = [double]"0.6088436993390358"
$ws.Range("S15").Value = [double]"-2.279160944206359"
$ws.Range("U15").Value = [double]"6.506873032118446"
$ws.Range("C16").Value = [double]"37.73812184455447"
$ws.Range("D16").Value = [double]"9.862893196220611"
$ws.Range("E16").Value = [double]"1"
$ws.Range("F16").Value = [double]"0.4931446598110306"
$ws.Range("G16").Value = [double]"-0.7916354698669097"
$ws.Range("J16").Value = [double]"0.2292656064641041"
$ws.Range("O16").Value = [double]"4.524986705778389E-06"
$ws.Range("P16").Value = [double]"9.862888671233906"
$ws.Range("T16").Value = [double]"-3.235244925562745"
$ws.Range("U16").Value = [double]"3.235244925562745"
$ws.Range("C17").Value = [double]"87.05258782565753"
$ws.Range("D17").Value = [double]"-4.184913538084566"
$ws.Range("E17").Value = [double]"1"
$ws.Range("F17").Value = [double]"0.2092456769042283"
$ws.Range("G17").Value = [double]"1.598709494803799"
$ws.Range("H17").Value = [double]"0.6965512767740842"
$ws.Range("J17").Value = [double]"0.7517985431907311"
$ws.Range("R17").Value = [double]"4.184913538084566"
$ws.Range("S17").Value = [double]"-2.920442868084405"
$ws.Range("U17").Value = [double]"4.516841769183306"
$ws.Range("C18").Value = [double]"51.52580579481268"
$ws.Range("D18").Value = [double]"-2.450665837139279"
$ws.Range("E18").Value = [double]"1"
$ws.Range("F18").Value = [double]"0.1225332918569639"
$ws.Range("G18").Value = [double]"0.9876412490563953"
$ws.Range("H18").Value = [double]"0.316016883925394"
$ws.Range("J18").Value = [double]"0.5605381194742565"
$ws.Range("R18").Value = [double]"2.450665837139279"
$ws.Range("S18").Value = [double]"-3.854495321823256"
$ws.Range("U18").Value = [double]"4.346893124021058"
$ws.Range("E19").Value = [double]"1"
$ws.Range("G19").Value = [double]"-1.598030624329664"
$ws.Range("J19").Value = [double]"0.002988844061765045"
$ws.Range("O19").Value = [double]"3.721591209718267E-07"
$ws.Range("P19").Value = [double]"9.570888225101037"
$ws.Range("T19").Value = [double]"-0.03471391355516474"
$ws.Range("U19").Value = [double]"0.03471391355516474"
$ws.Range("E20").Value = [double]"1"
$ws.Range("G20").Value = [double]"0.2180012696364086"
$ws.Range("J20").Value = [double]"0.3986845813368131"
$ws.Range("S20").Value = [double]"-4.00575795847613"
$ws.Range("U20").Value = [double]"4.00575795847613"
$ws.Range("C21").Value = [double]"47.82565319392014"
$ws.Range("D21").Value = [double]"-1.747547900335359"
$ws.Range("E21").Value = [double]"1"
$ws.Range("F21").Value = [double]"0.08737739501676793"
$ws.Range("G21").Value = [double]"0.3916281187605861"
$ws.Range("H21").Value = [double]"0.163572733738025"
$ws.Range("J21").Value = [double]"0.3573306601091225"
$ws.Range("R21").Value = [double]"1.747547900335359"
$ws.Range("S21").Value = [double]"-3.81758273844867"
$ws.Range("U21").Value = [double]"3.81758273844867"
$ws.Range("D22").Value = [double]"7.522354737325383"
$ws.Range("E22").Value = [double]"1"
$ws.Range("F22").Value = [double]"0.3761177368662691"
$ws.Range("G22").Value = [double]"-1.546343000520817"
$ws.Range("P22").Value = [double]"7.522354737325383"
$ws.Range("C23").Value = [double]"57.61177368662691"
$ws.Range("E23").Value = [double]"1"
$ws.Range("G23").Value = [double]"0.4647082704530749"
$ws.Range("J23").Value = [double]"0.4960081419217461"
$ws.Range("S23").Value = [double]"-3.916254150764301"
$ws.Range("U23").Value = [double]"3.916254150764301"
$ws.Range("C24").Value = [double]"38.03050293280541"
$ws.Range("D24").Value = [double]"-1.234185405829535E-05"
$ws.Range("E24").Value = [double]"1"
$ws.Range("F24").Value = [double]"6.170927029147677E-07"
$ws.Range("G24").Value = [double]"0.368377904220417"
$ws.Range("H24").Value = [double]"1.050643087625912E-06"
$ws.Range("J24").Value = [double]"0.3069807558714387"
$ws.Range("R24").Value = [double]"1.234185405829535E-05"
$ws.Range("S24").Value = [double]"-3.606088244707024"
$ws.Range("U24").Value = [double]"3.606088244707024"
$ws.Range("E25").Value = [double]"1"

Write-Output "Applied 164 cell updates"
